$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H62").Value = 3848.75
$ws.Range("I62").Value = 3700
$ws.Range("K62").Value = 3700
$ws.Range("M62").Value = -3076
$ws.Range("H65").Value = 3848.75
$ws.Range("I65").Value = 3700
$ws.Range("K65").Value = 18500
$ws.Range("M65").Value = -15380
$ws.Range("H76").Value = 5202.857
$ws.Range("I76").Value = 4403.3335
$ws.Range("J76").Value = 10000
$ws.Range("K76").Value = 4403.3335
$ws.Range("L76").Value = 10000
$ws.Range("M76").Value = -4088.3335
$ws.Range("N76").Value = -10630
$ws.Range("H79").Value = 5202.857
$ws.Range("I79").Value = 4403.3335
$ws.Range("J79").Value = 10000
$ws.Range("K79").Value = 4403.3335
$ws.Range("L79").Value = 10000
$ws.Range("M79").Value = -3311.3335
$ws.Range("N79").Value = -12184
$ws.Range("H86").Value = 3674.75
$ws.Range("I86").Value = 2900
$ws.Range("J86").Value = 4449.5
$ws.Range("K86").Value = 2900
$ws.Range("L86").Value = 4449.5
$ws.Range("M86").Value = -1777
$ws.Range("N86").Value = -6695.5
$ws.Range("H89").Value = 3674.75
$ws.Range("I89").Value = 2900
$ws.Range("J89").Value = 4449.5
$ws.Range("K89").Value = 14500
$ws.Range("L89").Value = 22247.5
$ws.Range("M89").Value = -8884
$ws.Range("N89").Value = -33479.5
$ws.Range("H92").Value = 2338.1365
$ws.Range("J92").Value = 3065.8
$ws.Range("L92").Value = 3065.8
$ws.Range("N92").Value = -5561.8
$ws.Range("H98").Value = 2533.4048
$ws.Range("I98").Value = 2321.1035
$ws.Range("K98").Value = 2321.1035
$ws.Range("M98").Value = -823.1035000000002
$ws.Range("H99").Value = 1752.25
$ws.Range("I99").Value = 323.5
$ws.Range("J99").Value = 3181
$ws.Range("K99").Value = 970.5
$ws.Range("L99").Value = 9543
$ws.Range("M99").Value = 527.5
$ws.Range("N99").Value = -12539
$ws.Range("H106").Value = 16645.334
$ws.Range("I106").Value = 18038.625
$ws.Range("K106").Value = 18038.625
$ws.Range("M106").Value = -17407.625
$ws.Range("H122").Value = 2533.4048
$ws.Range("I122").Value = 2321.1035
$ws.Range("K122").Value = 6963.310500000001
$ws.Range("M122").Value = -4513.310500000001
$ws.Range("H125").Value = 2833.2856
$ws.Range("I125").Value = 432.66666
$ws.Range("J125").Value = 4633.75
$ws.Range("K125").Value = 3893.99994
$ws.Range("L125").Value = 41703.75
$ws.Range("M125").Value = -1433.99994
$ws.Range("N125").Value = -46623.75
$ws.Range("H134").Value = 41747.5
$ws.Range("J134").Value = 41747.5
$ws.Range("L134").Value = 41747.5
$ws.Range("N134").Value = -51887.5
$ws.Range("H137").Value = 4075.9167
$ws.Range("J137").Value = 6744.3335
$ws.Range("L137").Value = 20233.0005
$ws.Range("N137").Value = -25333.0005
$ws.Range("H140").Value = 98999
$ws.Range("J140").Value = 98999
$ws.Range("L140").Value = 98999
$ws.Range("N140").Value = -109359

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H140").Value = 89922
$ws.Range("J140").Value = 89922
$ws.Range("L140").Value = 89922
$ws.Range("N140").Value = -100282

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H82").Value = 48988.6
$ws.Range("J82").Value = 74981
$ws.Range("L82").Value = 74981
$ws.Range("N82").Value = -75747
$ws.Range("H85").Value = 48988.6
$ws.Range("J85").Value = 74981
$ws.Range("L85").Value = 74981
$ws.Range("N85").Value = -77633

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 2578.853
$ws.Range("I31").Value = 1459
$ws.Range("K31").Value = 1459
$ws.Range("M31").Value = -1164
$ws.Range("H34").Value = 2578.853
$ws.Range("I34").Value = 1459
$ws.Range("K34").Value = 1459
$ws.Range("M34").Value = -1257
$ws.Range("H107").Value = 532.125
$ws.Range("I107").Value = 305.33334
$ws.Range("K107").Value = 305.33334
$ws.Range("M107").Value = 1614.66666
$ws.Range("H132").Value = 4598.8696
$ws.Range("I132").Value = 4451.263
$ws.Range("J132").Value = 5300
$ws.Range("K132").Value = 13353.789
$ws.Range("L132").Value = 15900
$ws.Range("M132").Value = -10823.789
$ws.Range("N132").Value = -20960
$ws.Range("H134").Value = 5728.091
$ws.Range("I134").Value = 5795.763
$ws.Range("J134").Value = 5299.5
$ws.Range("K134").Value = 17387.289
$ws.Range("L134").Value = 15898.5
$ws.Range("M134").Value = -14852.289
$ws.Range("N134").Value = -20968.5
$ws.Range("H141").Value = 280000
$ws.Range("J141").Value = 280000
$ws.Range("L141").Value = 280000
$ws.Range("N141").Value = -290360

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H34").Value = 6327.4287
$ws.Range("J34").Value = 6327.4287
$ws.Range("L34").Value = 18982.2861
$ws.Range("N34").Value = -19150.2861
$ws.Range("H39").Value = 10500.333
$ws.Range("J39").Value = 10500.333
$ws.Range("L39").Value = 31500.999
$ws.Range("N39").Value = -32088.999
$ws.Range("H55").Value = 10562.25
$ws.Range("J55").Value = 10562.25
$ws.Range("L55").Value = 31686.75
$ws.Range("N55").Value = -32040.75
$ws.Range("H68").Value = 2037.7858
$ws.Range("J68").Value = 2230.182
$ws.Range("L68").Value = 6690.545999999999
$ws.Range("N68").Value = -8312.545999999998
$ws.Range("H71").Value = 2037.7858
$ws.Range("J71").Value = 2230.182
$ws.Range("L71").Value = 20071.638
$ws.Range("N71").Value = -28183.638

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H126").Value = 6357.353
$ws.Range("J126").Value = 9710
$ws.Range("L126").Value = 29130
$ws.Range("N126").Value = -34070
$ws.Range("H132").Value = 7300.375
$ws.Range("I132").Value = 8400.5
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 25201.5
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -22671.5
$ws.Range("N132").Value = -17060
$ws.Range("H138").Value = 82999.92999999999
$ws.Range("J138").Value = 82999.92999999999
$ws.Range("L138").Value = 82999.92999999999
$ws.Range("N138").Value = -93279.92999999999

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 1286.4117
$ws.Range("I22").Value = 999
$ws.Range("J22").Value = 1609.75
$ws.Range("K22").Value = 999
$ws.Range("L22").Value = 1609.75
$ws.Range("M22").Value = -704
$ws.Range("N22").Value = -2199.75
$ws.Range("H27").Value = 1286.4117
$ws.Range("I27").Value = 999
$ws.Range("J27").Value = 1609.75
$ws.Range("K27").Value = 999
$ws.Range("L27").Value = 1609.75
$ws.Range("M27").Value = -892
$ws.Range("N27").Value = -1823.75
$ws.Range("H132").Value = 3972.1333
$ws.Range("I132").Value = 3166.3333
$ws.Range("K132").Value = 9498.999899999999
$ws.Range("M132").Value = -6968.999899999999

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H54").Value = 35998.7
$ws.Range("J54").Value = 49993.5
$ws.Range("L54").Value = 49993.5
$ws.Range("N54").Value = -51033.5
$ws.Range("H107").Value = 6284
$ws.Range("I107").Value = 9719.727999999999
$ws.Range("J107").Value = 885
$ws.Range("K107").Value = 29159.184
$ws.Range("L107").Value = 2655
$ws.Range("M107").Value = -27239.184
$ws.Range("N107").Value = -6495
$ws.Range("H126").Value = 1152.15
$ws.Range("J126").Value = 983.5
$ws.Range("L126").Value = 2950.5
$ws.Range("N126").Value = -7890.5
